$wb = $excel.ActiveWorkbook

$oldGuid = "d120e4fd-5ce2-49e9-aba3-09a2b49ec8aa"
$newGuid = "694e66a3-dd9d-40a0-babf-657fd9920659"

$oldXlfHash = "db2273b2c4a71febf59e71415699cafeaac3f6ce"
$newXlfHash = "6720df9041da850ddd9d787fc0a1c50402ef4d3a"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = ($newGuid + ".md")
$wsOverview.Range("B2").Value = ("e2e\" + $newGuid + ".md")
$wsOverview.Range("G2").Value = "2016-09-03 13:05:07"

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$2') {
        $h.TextToDisplay = ("e2e\" + $newGuid + ".md")
    }
}

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = ($newGuid + ".md")
$wsZh.Range("G2").Value = ($newGuid + "." + $newXlfHash + ".zh-cn.xlf")
$wsZh.Range("H2").Value = "2016-09-03 13:04:58"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = ($newGuid + ".md")
    }
}

foreach ($h in $wsZh.Hyperlinks) {
    if ($h.Range.Address() -eq '$I$2') {
        $h.Delete()
    }
}

$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""

$wsZh.Columns.Item(9).ColumnWidth = 17.835
$wsZh.Columns.Item(10).ColumnWidth = 20.835

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = ($newGuid + ".md")
$wsDe.Range("G2").Value = ($newGuid + "." + $newXlfHash + ".de-de.xlf")
$wsDe.Range("H2").Value = "2016-09-03 13:05:07"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = ($newGuid + ".md")
    }
}

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$I$2') {
        $h.Delete()
    }
}

$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""

$wsDe.Columns.Item(9).ColumnWidth = 17.835
$wsDe.Columns.Item(10).ColumnWidth = 20.835

Write-Output "edit complete"
